# MHV-17222: bump CPM CodeSystem metadata (version/date) and make sure the
# "wrap text / top align" formatting that was already present on the
# header and data rows is actually flagged as applied.

$wb = $excel.ActiveWorkbook

$wsMetadata = $wb.Worksheets.Item("Metadata")
$wsConcepts = $wb.Worksheets.Item("Concepts")

# --- bump Version + Date values on the Metadata sheet -----------------
$wsMetadata.Range("B3").Value = "0.2.10-beta"
$wsMetadata.Range("B8").Value = "2023-12-06T12:46:33-06:00"

# --- (re)apply the vertical-top / wrap-text alignment so it "sticks" --
# The cell styles already carried vertical="top" wrapText="true" in the
# alignment record, but the xf entries were missing applyAlignment, so
# Excel was not actually honoring it. Re-asserting WrapText / VerticalAlignment
# on the already-formatted ranges turns that flag on for both the header
# row style and the body-row style on each sheet.

$metaHeader = $wsMetadata.Range("A1:B1")
$metaHeader.VerticalAlignment = -4160
$metaHeader.WrapText = $true

$metaBody = $wsMetadata.Range("A2:B23")
$metaBody.VerticalAlignment = -4160
$metaBody.WrapText = $true

$conceptsHeader = $wsConcepts.Range("A1:D1")
$conceptsHeader.VerticalAlignment = -4160
$conceptsHeader.WrapText = $true

$conceptsBody = $wsConcepts.Range("A2:D2")
$conceptsBody.VerticalAlignment = -4160
$conceptsBody.WrapText = $true
